$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-62)
# from 45221 (2023-10-22) to 45224 (2023-10-25).
$ws.Range("C2:C62").Value = 45224
